# Generate Report for Handoff
# Replace the old run's GUID / content-hash / timestamps with the new run's
# values across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# New handoff run id (replaces 974531bb-c52e-4a1d-aceb-2541c48169a7)
$newGuid = "a8b4fae4-967f-472d-9957-90031e2f676a"

# New xliff content hash (replaces 6bb28891fbd352d4ebc9c146d5673786cc7cbdfc)
$newHash = "3b1ac6a24fb32e1289d9ab94fe8bc13a94c3686d"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-30 21:05:34"
foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = "e2e\$newGuid.md"
}

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-30 21:05:29"
foreach ($h in $wsZhCn.Hyperlinks) {
    $h.TextToDisplay = "$newGuid.md"
}

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-30 21:05:34"
foreach ($h in $wsDeDe.Hyperlinks) {
    $h.TextToDisplay = "$newGuid.md"
}
